# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels, reusing the existing bold/bordered
#     header style from the last pre-existing header cell (AC1) so the
#     new cells share style index "1" just like every other header cell. ---
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-45): same season record repeated for every player row. ---
$ws.Range("AD2:AD45").Value = 88
$ws.Range("AE2:AE45").Value = 74
$ws.Range("AF2:AF45").Value = 0

Write-Host "Season record columns added"
